# Applies the cryptos-list refresh described in the commit diff.
# Numeric-looking Price cells are prefixed with a literal apostrophe so Excel
# stores/keeps them as text (preserving trailing zeros / exact formatting)
# instead of silently re-parsing them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.955.15'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.415.73'
$ws.Range("E3").Value = '  +3.56%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''256.48'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").Value = '''661.23'
$ws.Range("E6").Value = '  +6.18%  '
$ws.Range("D7").Value = '''1.48'
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("D8").Value = '''0.435'
$ws.Range("E8").Value = '  +9.30%  '
$ws.Range("D9").Value = '''1.07'
$ws.Range("E9").Value = '  +11.70%  '
$ws.Range("D10").Value = '''0.998'
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("D11").Value = '3.410.74'
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("D12").Value = '''0.214'
$ws.Range("E12").Value = '  +7.06%  '
$ws.Range("D13").Value = '''41.97'
$ws.Range("E13").Value = '  +6.80%  '
$ws.Range("D14").Value = '''6.26'
$ws.Range("E14").Value = '  +14.47%  '
$ws.Range("D15").Value = '''0.0000260'
$ws.Range("E15").Value = '  +5.09%  '
$ws.Range("D16").Value = '97.396.37'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '4.046.18'
$ws.Range("E17").Value = '  +3.50%  '
$ws.Range("D18").Value = '''8.58'
$ws.Range("E18").Value = '  +36.33%  '
$ws.Range("D19").Value = '3.401.36'
$ws.Range("E19").Value = '  +3.21%  '
$ws.Range("D20").Value = '''17.47'
$ws.Range("E20").Value = '  +12.99%  '
$ws.Range("D21").Value = '''0.500'
$ws.Range("E21").Value = '  +66.44%  '
$ws.Range("D22").Value = '''10.86'
$ws.Range("E22").Value = '  +15.17%  '
$ws.Range("D23").Value = '''3.48'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '''511.15'
$ws.Range("E24").Value = '  +5.48%  '
$ws.Range("D25").Value = '''0.0000208'
$ws.Range("E25").Value = '  +3.06%  '
$ws.Range("D26").Value = '''6.26'
$ws.Range("E26").Value = '  +11.63%  '
$ws.Range("D27").Value = '''99.90'
$ws.Range("E27").Value = '  +12.63%  '
$ws.Range("D28").Value = '''12.62'
$ws.Range("E28").Value = '  +5.32%  '
$ws.Range("D29").Value = '3.585.99'
$ws.Range("E29").Value = '  +3.10%  '
$ws.Range("D30").Value = '''0.153'
$ws.Range("E30").Value = '  +12.26%  '
$ws.Range("D31").Value = '''11.60'
$ws.Range("E31").Value = '  +15.81%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = '''0.996'
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").Value = '''0.194'
$ws.Range("E33").Value = '  +3.59%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").Value = '''0.569'
$ws.Range("E35").Value = '  +22.97%  '
$ws.Range("D36").Value = '''29.82'
$ws.Range("E36").Value = '  +6.76%  '
$ws.Range("D37").Value = '''2.19'
$ws.Range("E37").Value = '  +12.97%  '
$ws.Range("D38").Value = '''7.81'
$ws.Range("E38").Value = '  +9.26%  '
$ws.Range("D39").Value = '''0.156'
$ws.Range("E39").Value = '  +6.08%  '
$ws.Range("D40").Value = '''517.98'
$ws.Range("E40").Value = '  +5.89%  '
$ws.Range("D41").Value = '''1.39'
$ws.Range("E41").Value = '  +13.01%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").Value = '''24.72'
$ws.Range("E42").Value = '  -0.41%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0437'
$ws.Range("E43").Value = '  +32.99%  '
$ws.Range("D44").Value = '''0.848'
$ws.Range("E44").Value = '  +7.77%  '
$ws.Range("D45").Value = '''3.71'
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").Value = '''3.34'
$ws.Range("E46").Value = '  +7.49%  '
$ws.Range("D47").Value = '''8.31'
$ws.Range("E47").Value = '  +13.57%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '''1.00'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = '''5.32'
$ws.Range("E49").Value = '  +12.75%  '
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").Value = '''1.59'
$ws.Range("E50").Value = '  +17.54%  '
$ws.Range("D51").Value = '''2.08'
$ws.Range("E51").Value = '  +7.74%  '
